$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven update of the cryptocurrency table to reflect the latest snapshot.
# - Columns B (Coin) and C (Link) hold plain text and are written directly.
# - Columns D (Price) and E (Volume(1h)) often look numeric (e.g. "1.00",
#   "0.0000246", "  -4.81%  ") which Excel would otherwise silently convert to
#   a number/date. Setting NumberFormat to "@" (Text) first keeps them as the
#   literal text strings found in the source data, matching the original file.

$rowUpdates = @(
    @{ Row = 2; D = "93.496.93"; E = "  -4.81%  " },
    @{ Row = 3; D = "3.387.89"; E = "  -0.66%  " },
    @{ Row = 4; E = "  +0.09%  " },
    @{ Row = 5; D = "234.65"; E = "  -7.79%  " },
    @{ Row = 6; D = "636.76"; E = "  -3.79%  " },
    @{ Row = 7; E = "  -4.90%  " },
    @{ Row = 8; D = "0.395"; E = "  -8.31%  " },
    @{ Row = 9; E = "  +0.16%  " },
    @{ Row = 10; D = "0.952"; E = "  -8.74%  " },
    @{ Row = 11; D = "3.382.70"; E = "  -0.72%  " },
    @{ Row = 12; E = "  -5.53%  " },
    @{ Row = 13; D = "41.00"; E = "  -9.07%  " },
    @{ Row = 14; D = "6.11"; E = "  -0.18%  " },
    @{ Row = 15; D = "93.578.47"; E = "  -4.55%  " },
    @{ Row = 16; D = "4.021.53"; E = "  -0.52%  " },
    @{ Row = 17; D = "0.0000246"; E = "  -4.91%  " },
    @{ Row = 18; D = "8.16"; E = "  -11.42%  " },
    @{ Row = 19; D = "3.385.10"; E = "  -1.62%  " },
    @{ Row = 20; D = "17.13"; E = "  -6.46%  " },
    @{ Row = 21; D = "11.23"; E = "  -2.45%  " },
    @{ Row = 22; D = "492.63"; E = "  -3.97%  " },
    @{ Row = 23; D = "0.464"; E = "  -11.22%  " },
    @{ Row = 24; D = "3.20"; E = "  -6.97%  " },
    @{ Row = 25; D = "0.0000188"; E = "  -6.98%  " },
    @{ Row = 26; D = "6.37"; E = "  -7.57%  " },
    @{ Row = 27; D = "90.47"; E = "  -7.09%  " },
    @{ Row = 28; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "11.67"; E = "  -6.50%  " },
    @{ Row = 29; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "11.38"; E = "  -6.86%  " },
    @{ Row = 30; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  +0.09%  " },
    @{ Row = 31; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "2.69"; E = "  -4.74%  " },
    @{ Row = 32; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.134"; E = "  -6.60%  " },
    @{ Row = 33; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "0.996"; E = "  -0.46%  " },
    @{ Row = 34; B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.174"; E = "  -7.96%  " },
    @{ Row = 35; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "28.90"; E = "  -1.09%  " },
    @{ Row = 36; B = "PolygonEcosystemToken"; C = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D = "0.538"; E = "  -5.14%  " },
    @{ Row = 37; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "534.62"; E = "  +1.60%  " },
    @{ Row = 38; B = "RenderToken"; C = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D = "7.53"; E = "  -6.08%  " },
    @{ Row = 39; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "1.42"; E = "  -5.16%  " },
    @{ Row = 40; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "1.00"; E = "  -0.05%  " },
    @{ Row = 41; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.148"; E = "  -3.46%  " },
    @{ Row = 42; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "0.891"; E = "  +3.12%  " },
    @{ Row = 43; B = "WhiteBITCoin"; C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D = "24.05"; E = "  -1.51%  " },
    @{ Row = 44; B = "MantraDAO"; C = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"; D = "3.64"; E = "  -1.24%  " },
    @{ Row = 45; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.68"; E = "  -3.47%  " },
    @{ Row = 46; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "5.59"; E = "  -1.02%  " },
    @{ Row = 47; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "2.15"; E = "  -5.10%  " },
    @{ Row = 48; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0399"; E = "  -6.72%  " },
    @{ Row = 49; D = "53.82"; E = "  -4.00%  " },
    @{ Row = 50; D = "3.23"; E = "  -0.74%  " },
    @{ Row = 51; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "7.94"; E = "  -7.98%  " }
)

foreach ($u in $rowUpdates) {
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $eCell = $ws.Cells.Item($u.Row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $u.E
    }
}

Write-Output ("Updated " + $rowUpdates.Count + " rows")
